$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 29898
$ws.Range("E2").Value = 580986334600
$ws.Range("F2").Value = 6567055351
$ws.Range("G2").Value = 0.08076

# Row 3
$ws.Range("D3").Value = 1875.36
$ws.Range("E3").Value = 225363437431
$ws.Range("F3").Value = 5168980366
$ws.Range("G3").Value = -0.64466

# Row 4
$ws.Range("D4").Value = 0.999833
$ws.Range("E4").Value = 83780327503
$ws.Range("F4").Value = 6062061027
$ws.Range("G4").Value = -0.04711

# Row 5
$ws.Range("D5").Value = 0.742252
$ws.Range("E5").Value = 38974495678
$ws.Range("F5").Value = 1572076246
$ws.Range("G5").Value = -3.75616

# Row 6
$ws.Range("D6").Value = 242.3
$ws.Range("E6").Value = 37734786692
$ws.Range("F6").Value = 220685302
$ws.Range("G6").Value = -0.16139

# Row 7
$ws.Range("D7").Value = 0.99993
$ws.Range("E7").Value = 26852220692
$ws.Range("F7").Value = 2754304258
$ws.Range("G7").Value = -0.01467

# Row 8
$ws.Range("D8").Value = 1874.18
$ws.Range("E8").Value = 14596700721
$ws.Range("F8").Value = 18994530
$ws.Range("G8").Value = -0.66086

# Row 9
$ws.Range("D9").Value = 0.314902
$ws.Range("E9").Value = 11032322205
$ws.Range("F9").Value = 159912431
$ws.Range("G9").Value = 1.02945

# Row 10
$ws.Range("B10").Value = "DOGE"
$ws.Range("C10").Value = "Dogecoin"
$ws.Range("D10").Value = 0.072154
$ws.Range("E10").Value = 10107512042
$ws.Range("F10").Value = 322702449
$ws.Range("G10").Value = 0.51235

# Row 11
$ws.Range("B11").Value = "SOL"
$ws.Range("C11").Value = "Solana"
$ws.Range("D11").Value = 24.69
$ws.Range("E11").Value = 9968187579
$ws.Range("F11").Value = 368966128
$ws.Range("G11").Value = -3.51955

# Row 12
$ws.Range("B12").Value = "TRX"
$ws.Range("C12").Value = "TRON"
$ws.Range("D12").Value = 0.083651
$ws.Range("E12").Value = 7509214563
$ws.Range("F12").Value = 324343535
$ws.Range("G12").Value = -2.7571

# Row 13
$ws.Range("B13").Value = "MATIC"
$ws.Range("C13").Value = "Polygon"
$ws.Range("D13").Value = 0.752761
$ws.Range("E13").Value = 7008942776
$ws.Range("F13").Value = 149024466
$ws.Range("G13").Value = -1.2753

# Row 14
$ws.Range("B14").Value = "DOT"
$ws.Range("C14").Value = "Polkadot"
$ws.Range("D14").Value = 5.42
$ws.Range("E14").Value = 6801666787
$ws.Range("F14").Value = 110138647
$ws.Range("G14").Value = 0.97168

# Row 15
$ws.Range("B15").Value = "LTC"
$ws.Range("C15").Value = "Litecoin"
$ws.Range("D15").Value = 92.55
$ws.Range("E15").Value = 6791480445
$ws.Range("F15").Value = 557117996
$ws.Range("G15").Value = -1.02416

# Row 16
$ws.Range("B16").Value = "BCH"
$ws.Range("C16").Value = "Bitcoin Cash"
$ws.Range("D16").Value = 252.28
$ws.Range("E16").Value = 4905983265
$ws.Range("F16").Value = 220550362
$ws.Range("G16").Value = 3.32662

# Row 17
$ws.Range("B17").Value = "WBTC"
$ws.Range("C17").Value = "Wrapped Bitcoin"
$ws.Range("D17").Value = 29875
$ws.Range("E17").Value = 4814361702
$ws.Range("F17").Value = 52848394
$ws.Range("G17").Value = -0.05935

# Row 18
$ws.Range("B18").Value = "AVAX"
$ws.Range("C18").Value = "Avalanche"
$ws.Range("D18").Value = 13.57
$ws.Range("E18").Value = 4694575001
$ws.Range("F18").Value = 98636737
$ws.Range("G18").Value = -1.37094

# Row 19
$ws.Range("D19").Value = 0.000007839999999999999475818607
$ws.Range("E19").Value = 4622799332
$ws.Range("F19").Value = 54224613
$ws.Range("G19").Value = 0.45531

# Row 20
$ws.Range("D20").Value = 6.08
$ws.Range("E20").Value = 4569999549
$ws.Range("F20").Value = 102999939
$ws.Range("G20").Value = -1.48416

# Row 21
$ws.Range("E21").Value = 4234992128
$ws.Range("F21").Value = 54746658
$ws.Range("G21").Value = 0.05196

# Row 22
$ws.Range("B22").Value = "XLM"
$ws.Range("C22").Value = "Stellar"
$ws.Range("D22").Value = 0.155064
$ws.Range("E22").Value = 4222602628
$ws.Range("F22").Value = 176012334
$ws.Range("G22").Value = -5.87573

# Row 23
$ws.Range("B23").Value = "LINK"
$ws.Range("C23").Value = "Chainlink"
$ws.Range("D23").Value = 8.02
$ws.Range("E23").Value = 4143991913
$ws.Range("F23").Value = 253962078
$ws.Range("G23").Value = 0.85036

# Row 24
$ws.Range("B24").Value = "BUSD"
$ws.Range("C24").Value = "Binance USD"
$ws.Range("D24").Value = 0.999604
$ws.Range("E24").Value = 3850889917
$ws.Range("F24").Value = 1441267226
$ws.Range("G24").Value = -0.02774

# Row 25
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "LEO Token"
$ws.Range("D25").Value = 4.02
$ws.Range("E25").Value = 3730118520
$ws.Range("F25").Value = 1967230
$ws.Range("G25").Value = -2.57592

# Row 26
$ws.Range("D26").Value = 164.96
$ws.Range("E26").Value = 2992147722
$ws.Range("F26").Value = 102390742
$ws.Range("G26").Value = 1.56219

# Row 27
$ws.Range("D27").Value = 0.999214
$ws.Range("E27").Value = 2804903993
$ws.Range("F27").Value = 983081241
$ws.Range("G27").Value = -0.00399

# Row 28
$ws.Range("D28").Value = 9.26
$ws.Range("E28").Value = 2706713855
$ws.Range("F28").Value = 55609814
$ws.Range("G28").Value = -1.10054

# Row 29
$ws.Range("D29").Value = 18.72
$ws.Range("E29").Value = 2659255878
$ws.Range("F29").Value = 93249591
$ws.Range("G29").Value = -0.05692

# Row 30
$ws.Range("D30").Value = 42.98
$ws.Range("E30").Value = 2578976709
$ws.Range("F30").Value = 1786469
$ws.Range("G30").Value = -0.01688

# Row 31
$ws.Range("B31").Value = "TON"
$ws.Range("C31").Value = "Toncoin"
$ws.Range("D31").Value = 1.51
$ws.Range("E31").Value = 2224389858
$ws.Range("F31").Value = 36244226
$ws.Range("G31").Value = 4.64268

# Row 32
$ws.Range("B32").Value = "FIL"
$ws.Range("C32").Value = "Filecoin"
$ws.Range("D32").Value = 4.61
$ws.Range("E32").Value = 2014581584
$ws.Range("F32").Value = 136856611
$ws.Range("G32").Value = 2.22632

# Row 33
$ws.Range("B33").Value = "ICP"
$ws.Range("C33").Value = "Internet Computer"
$ws.Range("D33").Value = 4.31
$ws.Range("E33").Value = 1886624054
$ws.Range("F33").Value = 30943276
$ws.Range("G33").Value = 5.08884

# Row 34
$ws.Range("B34").Value = "LDO"
$ws.Range("C34").Value = "Lido DAO"
$ws.Range("D34").Value = 2.03
$ws.Range("E34").Value = 1784162289
$ws.Range("F34").Value = 30249874
$ws.Range("G34").Value = -0.38791

# Row 35
$ws.Range("B35").Value = "HBAR"
$ws.Range("C35").Value = "Hedera"
$ws.Range("D35").Value = 0.05333
$ws.Range("E35").Value = 1724584728
$ws.Range("F35").Value = 43330230
$ws.Range("G35").Value = -1.73472

# Row 36
$ws.Range("B36").Value = "APT"
$ws.Range("C36").Value = "Aptos"
$ws.Range("D36").Value = 7.62
$ws.Range("E36").Value = 1650453199
$ws.Range("F36").Value = 57840880
$ws.Range("G36").Value = -0.16152

# Row 37
$ws.Range("B37").Value = "CRO"
$ws.Range("C37").Value = "Cronos"
$ws.Range("D37").Value = 0.060651
$ws.Range("E37").Value = 1586543952
$ws.Range("F37").Value = 5511953
$ws.Range("G37").Value = -0.193

# Row 38
$ws.Range("B38").Value = "ARB"
$ws.Range("C38").Value = "Arbitrum"
$ws.Range("D38").Value = 1.24
$ws.Range("E38").Value = 1575608466
$ws.Range("F38").Value = 89718095
$ws.Range("G38").Value = -0.03452

# Row 39
$ws.Range("D39").Value = 103.58
$ws.Range("E39").Value = 1507102098
$ws.Range("F39").Value = 13883033
$ws.Range("G39").Value = 1.09311

# Row 40
$ws.Range("D40").Value = 0.01964414
$ws.Range("E40").Value = 1426896190
$ws.Range("F40").Value = 34625931
$ws.Range("G40").Value = 0.21474

# Row 41
$ws.Range("D41").Value = 1.46
$ws.Range("E41").Value = 1367080915
$ws.Range("F41").Value = 39438315
$ws.Range("G41").Value = -0.49028

# Row 42
$ws.Range("B42").Value = "OP"
$ws.Range("C42").Value = "Optimism"
$ws.Range("D42").Value = 1.57
$ws.Range("E42").Value = 1067536026
$ws.Range("F42").Value = 89711411
$ws.Range("G42").Value = 2.95097

# Row 43
$ws.Range("D43").Value = 0.116991
$ws.Range("E43").Value = 1064787770
$ws.Range("F43").Value = 18310684
$ws.Range("G43").Value = 0.2467

# Row 44
$ws.Range("B44").Value = "AAVE"
$ws.Range("C44").Value = "Aave"
$ws.Range("D44").Value = 72.37
$ws.Range("E44").Value = 1047921093
$ws.Range("F44").Value = 40156166
$ws.Range("G44").Value = -0.64191

# Row 45
$ws.Range("B45").Value = "MKR"
$ws.Range("C45").Value = "Maker"
$ws.Range("D45").Value = 1109.98
$ws.Range("E45").Value = 1000889234
$ws.Range("F45").Value = 62102236
$ws.Range("G45").Value = 0.07396

# Row 46
$ws.Range("B46").Value = "RETH"
$ws.Range("C46").Value = "Rocket Pool ETH"
$ws.Range("D46").Value = 2021.97
$ws.Range("E46").Value = 994839764
$ws.Range("F46").Value = 4503824
$ws.Range("G46").Value = -0.60471

# Row 47
$ws.Range("B47").Value = "SNX"
$ws.Range("C47").Value = "Synthetix Network"
$ws.Range("D47").Value = 2.91
$ws.Range("E47").Value = 932745423
$ws.Range("F47").Value = 69020758
$ws.Range("G47").Value = -2.1377

# Row 48
$ws.Range("B48").Value = "FRAX"
$ws.Range("C48").Value = "Frax"
$ws.Range("D48").Value = 0.998453
$ws.Range("E48").Value = 921050953
$ws.Range("F48").Value = 4144615
$ws.Range("G48").Value = 0.01445

# Row 49
$ws.Range("B49").Value = "ALGO"
$ws.Range("C49").Value = "Algorand"
$ws.Range("D49").Value = 0.115529
$ws.Range("E49").Value = 899761348
$ws.Range("F49").Value = 28532752
$ws.Range("G49").Value = -0.50173

# Row 50
$ws.Range("B50").Value = "SAND"
$ws.Range("C50").Value = "The Sandbox"
$ws.Range("D50").Value = 0.4543
$ws.Range("E50").Value = 879892274
$ws.Range("F50").Value = 47492284
$ws.Range("G50").Value = 1.86553

# Row 51
$ws.Range("B51").Value = "EGLD"
$ws.Range("C51").Value = "MultiversX"
$ws.Range("D51").Value = 33.88
$ws.Range("E51").Value = 869873872
$ws.Range("F51").Value = 7763618
$ws.Range("G51").Value = 0.21955
